$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as TEXT even if it looks numeric,
# to match the source data's string formatting (e.g. '249.55').
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$ws.Range("D2").Value = "36.656.53"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").Value = "1.922.47"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "249.55"
$ws.Range("E5").Value = "  +2.03%  "
Set-TextValue $ws.Range("D6") "0.697"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  -0.02%  "
Set-TextValue $ws.Range("D8") "44.39"
$ws.Range("E8").Value = "  +2.39%  "
Set-TextValue $ws.Range("D9") "58.57"
$ws.Range("E9").Value = "  +9.65%  "
Set-TextValue $ws.Range("D10") "0.367"
$ws.Range("E10").Value = "  +3.48%  "
Set-TextValue $ws.Range("D11") "0.0766"
$ws.Range("E11").Value = "  +3.54%  "
Set-TextValue $ws.Range("D12") "0.0997"
$ws.Range("E12").Value = "  +2.57%  "
Set-TextValue $ws.Range("D13") "14.53"
$ws.Range("E13").Value = "  +8.43%  "
Set-TextValue $ws.Range("D14") "0.798"
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("D15").Value = "2.202.00"
$ws.Range("E15").Value = "  +1.95%  "
Set-TextValue $ws.Range("D16") "5.13"
$ws.Range("E16").Value = "  +4.87%  "
$ws.Range("D17").Value = "1.920.50"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "36.593.23"
$ws.Range("E18").Value = "  +3.07%  "
Set-TextValue $ws.Range("D19") "74.48"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").Value = "0.0₃0860"
$ws.Range("E20").Value = "  +4.86%  "
Set-TextValue $ws.Range("D21") "250.95"
$ws.Range("E21").Value = "  +2.90%  "
Set-TextValue $ws.Range("D22") "13.27"
$ws.Range("E22").Value = "  +3.93%  "
Set-TextValue $ws.Range("D23") "5.21"
$ws.Range("E23").Value = "  +5.39%  "
Set-TextValue $ws.Range("D24") "2.70"
$ws.Range("E24").Value = "  +2.29%  "
Set-TextValue $ws.Range("D25") "0.999"
$ws.Range("E25").Value = "  -0.16%  "
Set-TextValue $ws.Range("D26") "2.21"
$ws.Range("E26").Value = "  +2.58%  "
Set-TextValue $ws.Range("D27") "168.08"
$ws.Range("E27").Value = "  +1.80%  "
Set-TextValue $ws.Range("D28") "8.86"
$ws.Range("E28").Value = "  +3.92%  "
Set-TextValue $ws.Range("D29") "18.81"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("E30").Value = "  +1.86%  "
Set-TextValue $ws.Range("D31") "4.55"
$ws.Range("E31").Value = "  +6.70%  "
Set-TextValue $ws.Range("D32") "0.0610"
$ws.Range("E32").Value = "  +3.83%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("E34").Value = "  +4.94%  "
$ws.Range("E35").Value = "  -0.01%  "
Set-TextValue $ws.Range("D36") "0.0847"
$ws.Range("E36").Value = "  +19.02%  "
$ws.Range("E37").Value = "  -13.19%  "
$ws.Range("E40").Value = "  +3.54%  "
Set-TextValue $ws.Range("D41") "108.07"
$ws.Range("E41").Value = "  +12.68%  "
$ws.Range("E42").Value = "  +5.23%  "
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("D45").Value = "1.342.14"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("E47").Value = "  +5.30%  "
Set-TextValue $ws.Range("D48") "0.0813"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("D51").Value = "2.099.66"
$ws.Range("E51").Value = "  +1.67%  "

# Row 38/39: Gas and ImmutableX swap positions (with updated values)
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D38") "0.875"
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("B39").Value = "Gas"
$ws.Range("C39").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
Set-TextValue $ws.Range("D39") "17.95"
$ws.Range("E39").Value = "  +47.54%  "
